# Add 2022 mortality data (new column S) to the Supp_Table_1 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: S1 = "2022" (stored as text, like the other year headers B1:R1)
$ws.Range("S1").Value = "'2022"

# New 2022 values for each ICD10 group row (rows 2-10)
$ws.Range("S2").Value = 21
$ws.Range("S3").Value = 19
$ws.Range("S4").Value = 255
$ws.Range("S5").Value = 34
$ws.Range("S6").Value = 19
$ws.Range("S7").Value = 136
$ws.Range("S8").Value = 95
$ws.Range("S9").Value = 7
$ws.Range("S10").Value = 120
